$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regression update: bump the PERIODE_BULANAN (Q2) and VERIFIKASI (T2)
# test-data values from the 2023 period to the 2024 period.
$ws.Range("Q2").Value = "202405"
$ws.Range("T2").Value = "15/04/2024"

# Reflect the author's new cursor position / selection when they saved.
$ws.Range("Q3").Select()
